$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing header cells C1:M1 right by one column into D1:N1, making room for
# the new "Units" column at C1. Column width metadata (the <cols> block) is left
# untouched, matching a direct-value shift rather than a true sheet/column insert.
# Read all the old values up front so writes don't clobber cells still to be read.
$oldValues = @()
for ($col = 3; $col -le 13; $col++) {
    $oldValues += ,$ws.Cells.Item(1, $col).Value2
}
for ($i = $oldValues.Length - 1; $i -ge 0; $i--) {
    $ws.Cells.Item(1, $i + 4).Value2 = $oldValues[$i]
}

# New header cell - copy formatting from the other header cells (e.g. B1) and set text.
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value2 = "Units"

# The shift above moved the old M1 ("Hi_Hi Limit") value into the previously-empty
# N1 cell, which doesn't carry the header formatting by default - apply it explicitly.
$ws.Range("B1").Copy()
$ws.Range("N1").PasteSpecial(-4122)

# Move the active selection to E5, matching the post-edit workbook state.
$ws.Range("E5").Select()
